$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.190.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.925.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +1.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.924.20"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.443"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000227"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.414.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.149.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.928.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "431.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("E22").Value = "  +2.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.32%  "
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.38%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("E34").Value = "  +2.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0869"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.18%  "
$ws.Range("E36").Value = "  +0.83%  "
$ws.Range("E37").Value = "  +3.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.289"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "380.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.71%  "
$ws.Range("E46").Value = "  +1.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.709.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.66%  "
$ws.Range("E51").Value = "  +0.59%  "
